$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# "u matrix: electricity use from gas boiler put to zero"
# Row 4 = Electricity "Need" row.
# H4 = Exploiting Gas boiler for Heating
# K4 = Exploiting Gas boiler for Hot Sanitary Water
# P4 = Exploiting Gas Stove for Cooking
$ws.Range("H4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("P4").Value = 0
